$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.340.25"
$ws.Range("E2").Value = "  -0.93%  "

$ws.Range("D3").Value = "'1.548.59"
$ws.Range("E3").Value = "  -1.97%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'209.77"
$ws.Range("E5").Value = "  -1.72%  "

$ws.Range("D6").Value = "'0.481"
$ws.Range("E6").Value = "  -1.57%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("D8").Value = "'23.87"
$ws.Range("E8").Value = "  -0.89%  "

$ws.Range("E9").Value = "  -1.90%  "

$ws.Range("E10").Value = "  -1.55%  "

$ws.Range("D11").Value = "'0.0887"
$ws.Range("E11").Value = "  -0.58%  "

$ws.Range("D12").Value = "'1.767.94"
$ws.Range("E12").Value = "  -2.05%  "

$ws.Range("D13").Value = "'1.556.08"
$ws.Range("E13").Value = "  -1.52%  "

$ws.Range("D14").Value = "'28.304.01"
$ws.Range("E14").Value = "  -0.99%  "

$ws.Range("D15").Value = "'3.62"
$ws.Range("E15").Value = "  -2.30%  "

$ws.Range("E16").Value = "  -2.00%  "

$ws.Range("D17").Value = "'60.86"
$ws.Range("E17").Value = "  -2.12%  "

$ws.Range("D18").Value = "'227.83"
$ws.Range("E18").Value = "  -1.51%  "

$ws.Range("D20").Value = "'0.0₃0676"
$ws.Range("E20").Value = "  -2.25%  "

$ws.Range("E21").Value = "  -0.04%  "

$ws.Range("E22").Value = "  -0.24%  "

$ws.Range("D23").Value = "'8.91"
$ws.Range("E23").Value = "  -3.00%  "

$ws.Range("E24").Value = "  -2.48%  "

$ws.Range("D25").Value = "'151.47"
$ws.Range("E25").Value = "  +0.23%  "

$ws.Range("D26").Value = "'14.72"
$ws.Range("E26").Value = "  -2.07%  "

$ws.Range("E27").Value = "  -1.11%  "

$ws.Range("E28").Value = "  -0.15%  "

$ws.Range("D29").Value = "'6.23"
$ws.Range("E29").Value = "  -3.30%  "

$ws.Range("D30").Value = "'0.0467"
$ws.Range("E30").Value = "  -3.62%  "

$ws.Range("E31").Value = "  -4.69%  "

$ws.Range("E32").Value = "  -1.35%  "

$ws.Range("D33").Value = "'1.384.51"
$ws.Range("E33").Value = "  -0.96%  "

$ws.Range("E34").Value = "  -3.37%  "

$ws.Range("D35").Value = "'1.07"
$ws.Range("E35").Value = "  +0.72%  "

$ws.Range("D36").Value = "'1.48"
$ws.Range("E36").Value = "  -3.63%  "

$ws.Range("D37").Value = "'2.33"
$ws.Range("E37").Value = "  -1.22%  "

$ws.Range("E38").Value = "  -2.70%  "

$ws.Range("E39").Value = "  -2.81%  "

$ws.Range("D40").Value = "'1.92"
$ws.Range("E40").Value = "  +2.02%  "

$ws.Range("D41").Value = "'0.508"
$ws.Range("E41").Value = "  -2.68%  "

$ws.Range("E42").Value = "  -0.15%  "

$ws.Range("D43").Value = "'0.773"
$ws.Range("E43").Value = "  -2.39%  "

$ws.Range("D44").Value = "'0.0453"
$ws.Range("E44").Value = "  -2.95%  "

$ws.Range("D45").Value = "'5.36"
$ws.Range("E45").Value = "  -1.73%  "

$ws.Range("D46").Value = "'61.89"
$ws.Range("E46").Value = "  -2.26%  "

$ws.Range("D47").Value = "'1.681.83"
$ws.Range("E47").Value = "  -2.02%  "

$ws.Range("D48").Value = "'0.864"
$ws.Range("E48").Value = "  -10.02%  "

$ws.Range("D49").Value = "'85.65"
$ws.Range("E49").Value = "  -0.98%  "

$ws.Range("D50").Value = "'42.03"
$ws.Range("E50").Value = "  +4.30%  "

$ws.Range("D51").Value = "'0.0₆0101"
$ws.Range("E51").Value = "  -1.63%  "
